$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold / bordered / centered) onto the two
# new header cells before setting their text, then fill in the labels.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 (col I) and IF (col J) columns, rows 2-40
$values = @(
    @(5,5),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(7,7),
    @(8,9),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(6,6),
    @(7,7),
    @(9,9),
    @(8,8),
    @(10,10),
    @(6,6),
    @(8,8),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(6,7),
    @(6,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(10,10),
    @(7,7),
    @(7,8),
    @(7,7),
    @(6,6),
    @(6,7),
    @(4,4),
    @(4,4),
    @(4,4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
